# Updates the "Estado de Cuenta" base data table (rows 16-22, NIT-9013218651
# sheet): the Periodo Mora / Valor Mora rows are re-sorted into ascending
# period order (2211, 2212, 2301, 2302, 2303, 2304, 2305), carrying each
# row's "Valor Mora" amount along with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending order for "Periodo Mora" (column E, rows 16-22)
$ws.Range("E16").Value = "2211"
$ws.Range("E17").Value = "2212"
$ws.Range("E18").Value = "2301"
$ws.Range("E19").Value = "2302"
$ws.Range("E20").Value = "2303"
$ws.Range("E21").Value = "2304"
$ws.Range("E22").Value = "2305"

# Matching "Valor Mora" (column F) values, carried along with their period
$ws.Range("F16").Value = 40000
$ws.Range("F17").Value = 40000
$ws.Range("F18").Value = 40000
$ws.Range("F19").Value = 40000
$ws.Range("F20").Value = 40000
$ws.Range("F21").Value = 40000
$ws.Range("F22").Value = 32000
